$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 6, shifting rows 6:25 down to 7:26
$ws.Rows.Item(6).Insert()

# Fill the new row 6 with data (most columns copy the constant values used throughout the sheet)
$ws.Cells.Item(6, 1).Value = 10
$ws.Cells.Item(6, 2).Value = 'Vega Modelo de Temuco'
$ws.Cells.Item(6, 3).Value = 'La Araucanía'
$ws.Cells.Item(6, 4).Value = 44525
$ws.Cells.Item(6, 5).Value = 9
$ws.Cells.Item(6, 6).Value = 'Fruta'
$ws.Cells.Item(6, 7).Value = 100103
$ws.Cells.Item(6, 8).Value = 'Frutos de hueso (carozo)'
$ws.Cells.Item(6, 9).Value = 100103003
$ws.Cells.Item(6, 10).Value = 'Damasco'
$ws.Cells.Item(6, 11).Value = 'Castle Brite'
$ws.Cells.Item(6, 12).Value = 'Primera'
$ws.Cells.Item(6, 13).Value = 55
$ws.Cells.Item(6, 14).Value = 20000
$ws.Cells.Item(6, 15).Value = 20000
$ws.Cells.Item(6, 16).Value = 20000
$ws.Cells.Item(6, 17).Value = '$/bandeja 10 kilos'
$ws.Cells.Item(6, 18).Value = 'Provincia de San Felipe de Aconcagua'
$ws.Cells.Item(6, 19).Value = 2000
$ws.Cells.Item(6, 20).Value = 10

# Match the date number format used by the other rows in column D
$ws.Cells.Item(6, 4).NumberFormat = $ws.Cells.Item(7, 4).NumberFormat
